# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the "05482384-..." file across the Overview, zh-cn and
# de-de sheets, and widens the "Error Detail" column so the message is
# readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# All four cells below originally shared the single "Ready for handoff"
# string; that shared text now reads "Handback transform failed", so
# every cell that displayed it must be updated in lockstep.

# Overview sheet: row 3 corresponds to 05482384-dec5-4d6a-8bce-6916aab1c32e.md
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# zh-cn sheet: row 3 = 05482384-... entry
$wsZhCn.Range("C3").Value = "Handback transform failed"
# Column P = "Error Detail"
$wsZhCn.Range("P3").Value = "Handback file name: 1rutsxsu.p2d is different with handoff file name: 05482384-dec5-4d6a-8bce-6916aab1c32e.df01bf7182275bc5d2252422b203a6bcd5f6af1b.zh-cn."
# ColumnWidth round-trips through a pixel-based conversion that adds
# 5/6 of a character to whatever is assigned, so back that off here to
# land on the target stored width of 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet: row 3 = 05482384-... entry
$wsDeDe.Range("C3").Value = "Handback transform failed"
# Column P = "Error Detail"
$wsDeDe.Range("P3").Value = "Handback file name: 1rutsxsu.p2d is different with handoff file name: 05482384-dec5-4d6a-8bce-6916aab1c32e.df01bf7182275bc5d2252422b203a6bcd5f6af1b.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
